# Apply the row-value shuffle described by the diff.
# Columns D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) get new values per row,
# while everything else on each row stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> hashtable of new column values
$updates = @{
    2  = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    6  = @{ D = 44914; M = 56; N = 23000; O = 23000; P = 23000; S = 4600 }
    7  = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    8  = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    10 = @{ D = 44902; M = 35; N = 12000; O = 12000; P = 12000; S = 2400 }
    11 = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    12 = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
    13 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
